$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.325.64"
$ws.Range("E2").Value = "  -1.40%  "

$ws.Range("D3").Value = "2.212.71"
$ws.Range("E3").Value = "  -4.96%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "295.60"
$ws.Range("E5").Value = "  -4.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "81.39"
$ws.Range("E6").Value = "  -3.86%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.511"
$ws.Range("E7").Value = "  -3.73%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.469"
$ws.Range("E9").Value = "  -3.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0771"
$ws.Range("E10").Value = "  -5.04%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "29.49"
$ws.Range("E11").Value = "  -1.80%  "

$ws.Range("E12").Value = "  -11.16%  "

$ws.Range("E13").Value = "  -2.71%  "

$ws.Range("D14").Value = "2.537.93"
$ws.Range("E14").Value = "  -5.41%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.24"
$ws.Range("E15").Value = "  -2.52%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.96"
$ws.Range("E16").Value = "  -5.12%  "

$ws.Range("D17").Value = "2.195.69"
$ws.Range("E17").Value = "  -4.76%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.712"
$ws.Range("E18").Value = "  -5.47%  "

$ws.Range("D19").Value = "39.219.35"
$ws.Range("E19").Value = "  -1.57%  "

$ws.Range("D20").Value = "0.0₃0872"
$ws.Range("E20").Value = "  -3.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.72"
$ws.Range("E21").Value = "  -5.95%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "64.63"
$ws.Range("E22").Value = "  -4.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.27"
$ws.Range("E23").Value = "  -3.32%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "228.22"
$ws.Range("E24").Value = "  -2.98%  "

$ws.Range("E25").Value = "  -0.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.40"
$ws.Range("E26").Value = "  -5.81%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.80"
$ws.Range("E27").Value = "  -0.13%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.59"
$ws.Range("E28").Value = "  -3.20%  "

$ws.Range("E29").Value = "  -1.21%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.08"
$ws.Range("E30").Value = "  -2.06%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "149.51"
$ws.Range("E31").Value = "  -1.87%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.70"
$ws.Range("E32").Value = "  -9.85%  "

$ws.Range("E33").Value = "  -0.09%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.80"
$ws.Range("E34").Value = "  -5.95%  "

$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0696"
$ws.Range("E35").Value = "  -3.36%  "

$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.34"
$ws.Range("E36").Value = "  -4.44%  "

$ws.Range("E37").Value = "  -3.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "15.48"
$ws.Range("E38").Value = "  -1.20%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0960"
$ws.Range("E39").Value = "  -3.58%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.64"
$ws.Range("E40").Value = "  -3.93%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.65"
$ws.Range("E41").Value = "  -3.20%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.61"
$ws.Range("E42").Value = "  -5.40%  "

$ws.Range("D43").Value = "1.903.10"
$ws.Range("E43").Value = "  -2.26%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.06"
$ws.Range("E44").Value = "  -8.97%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0259"
$ws.Range("E45").Value = "  -2.45%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.06"
$ws.Range("E46").Value = "  -2.26%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.10"
$ws.Range("E47").Value = "  -7.88%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.61"
$ws.Range("E48").Value = "  -2.28%  "

$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.407.50"
$ws.Range("E49").Value = "  -5.73%  "

$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.05"
$ws.Range("E50").Value = "  +0.82%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "87.63"
$ws.Range("E51").Value = "  -5.49%  "
